$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete all data rows (rows 2-5), keeping only the header row
$ws.Rows("2:5").Delete()

# Add the new column header in J1
$ws.Range("J1").Value = "custom_fields_json"
